$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "0ee1680a-1225-4fec-9439-5db9edb50fdf"
$ws.Range("A3").Value = "cb9b06b9-f483-4c37-a80c-8f8d4c9b2318"
$ws.Range("A4").Value = "d92218a9-cab0-45c5-b0c4-9144ecf40990"
$ws.Range("A5").Value = "e0d41b4d-880b-4d10-9ffd-dcab73d61046"
$ws.Range("A6").Value = "72da5ef6-a341-401d-9734-697b4c5e7ec9"
$ws.Range("A7").Value = "f64eef5d-1492-4f14-93df-8e5cf9cf8077"
$ws.Range("A8").Value = "c7cd5991-b261-4f6d-bc19-2c6c8414f4f6"
$ws.Range("A9").Value = "b6ffc9f8-f392-41fb-ba47-a30230394c22"
$ws.Range("A10").Value = "5744cbf7-95bc-4ed6-bb2e-b93312c91acd"
$ws.Range("A11").Value = "d925de4b-5427-4ed4-8e01-1f6e5ff83bb6"
$ws.Range("A12").Value = "1191073f-bcaa-4e80-81b9-4c9a3640e780"
$ws.Range("A13").Value = "7f3504f3-611f-4bd2-838f-169994668795"
$ws.Range("A14").Value = "eef165cb-8009-42f1-9669-574ff4900540"
$ws.Range("A15").Value = "671c670e-2ba9-4905-9aae-690cba81c085"
$ws.Range("A16").Value = "365ebfcf-b03a-40f6-a886-1c61b72b09b7"
$ws.Range("A17").Value = "9c378892-e1f5-4f19-9f47-c07029f0b5ab"
$ws.Range("A18").Value = "5b5f39d2-ce93-485b-87f4-7f3e7955dcfd"
$ws.Range("A19").Value = "0ccda549-b12a-4443-87ac-a2630657ae7d"
$ws.Range("A20").Value = "54b5fc24-edb7-4478-90f6-2c8a9cfe83c0"
$ws.Range("A21").Value = "82e857ab-3481-4542-860c-26a93d63bb7e"
$ws.Range("A22").Value = "69958f31-2426-4577-977b-65acb653967e"
$ws.Range("A23").Value = "f203b025-ec2d-4365-8b39-7bd68f0b322f"
$ws.Range("A24").Value = "f9e8ff5e-d73f-438f-9982-78a4cad27a6b"
$ws.Range("A25").Value = "f1ed3f6b-3573-413a-b8a1-b8d4b541fa62"
$ws.Range("A26").Value = "dd3996a1-ac4b-4f2f-9824-7d9540986655"
$ws.Range("A27").Value = "7c437e11-b031-49a1-a9fb-5d9ea101e50b"
$ws.Range("A28").Value = "d952d1f8-4cb4-41b9-8508-1e8d3d8a6263"
$ws.Range("A29").Value = "14a65320-cdb1-4637-a359-ebe1201e0f5a"
$ws.Range("A30").Value = "c120927c-23fb-4f7d-b954-06778eaaac70"
$ws.Range("A31").Value = "87521189-6ec1-42b7-bfa7-2aab4a458e01"
$ws.Range("A32").Value = "6b6efc2a-603d-4e56-9e37-9a7e459d3157"
$ws.Range("A33").Value = "316dde02-3de4-48f6-8de5-edc31f99ce07"
$ws.Range("A34").Value = "c1c52f42-14c2-492b-b09b-dbcc854a7d39"
$ws.Range("A35").Value = "e15cfed8-2be4-4964-8d07-586e2d3c0545"
$ws.Range("A36").Value = "4998b4bd-a968-4ce6-a019-8f07346259f0"
$ws.Range("A37").Value = "8a28926b-dca3-44c9-84f8-a3311aab2ca6"
$ws.Range("A38").Value = "f1de832b-df3a-465a-aae3-16376dd4d57e"
$ws.Range("A39").Value = "255e9c1e-0781-4f99-a533-54f7df356181"
$ws.Range("A40").Value = "4e5a17f8-06c4-4f98-83ed-c074bfe99235"
$ws.Range("A41").Value = "16cf7b2f-51df-4f64-a871-76a2a1ec4d87"
$ws.Range("A42").Value = "95c64599-a475-45bd-b10e-81aa0ca774df"
$ws.Range("A43").Value = "91fa2adb-ff12-4ac8-bae2-74dc500fd47f"
$ws.Range("A44").Value = "67caa053-810e-4bfd-a782-2f17f26345e0"
$ws.Range("A45").Value = "f1fc1d0e-1899-4fa6-881a-ad9079605158"
$ws.Range("A46").Value = "5d7a0e60-3797-468a-ada4-6e7eea1432f7"
$ws.Range("A47").Value = "94d1cbca-1bed-44d1-ade0-4178cc819fc7"
$ws.Range("A48").Value = "1a244aa8-caee-46ad-ba49-ba57937f5bf8"
$ws.Range("A49").Value = "b41293da-d4d7-481c-ad14-14de5045391c"
